$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.225969755228931
$ws.Range("D2").Value = 0.006298024303401206
$ws.Range("E2").Value = 0.1393848103014008
$ws.Range("F2").Value = 0.7533793398548454
$ws.Range("G2").Value = 0.6201500300062435
$ws.Range("H2").Value = 0.6454766957521656
$ws.Range("I2").Value = 0.753837131855704
$ws.Range("L2").Value = 0.1684885293038576
$ws.Range("O2").Value = 2.51363540873939
$ws.Range("C3").Value = 0.2233335232329381
$ws.Range("D3").Value = 0.006197097785202388
$ws.Range("E3").Value = 0.136342067306547
$ws.Range("F3").Value = 0.7228535965416398
$ws.Range("G3").Value = 0.5877630113455865
$ws.Range("H3").Value = 0.6345060536392282
$ws.Range("I3").Value = 0.7294046322488583
$ws.Range("L3").Value = 0.1634208008851488
$ws.Range("O3").Value = 2.42162939103568
$ws.Range("C4").Value = 0.2218256378161954
$ws.Range("D4").Value = 0.006133699267312664
$ws.Range("E4").Value = 0.1345519741413064
$ws.Range("F4").Value = 0.7045593754398425
$ws.Range("G4").Value = 0.5682508425061883
$ws.Range("H4").Value = 0.6281106818067173
$ws.Range("I4").Value = 0.7148438111719955
$ws.Range("L4").Value = 0.16040954894941
$ws.Range("O4").Value = 2.366627523460778
$ws.Range("C5").Value = 0.2212390354097948
$ws.Range("D5").Value = 0.006107505335622321
$ws.Range("E5").Value = 0.1338421385254591
$ws.Range("F5").Value = 0.6972170454625797
$ws.Range("G5").Value = 0.5603932113413919
$ws.Range("H5").Value = 0.6255902006696203
$ws.Range("I5").Value = 0.709020924664209
$ws.Range("L5").Value = 0.1592076312723947
$ws.Range("O5").Value = 2.34458814483088
$ws.Range("C6").Value = 0.2211433145577928
$ws.Range("D6").Value = 0.006103134225755369
$ws.Range("E6").Value = 0.1337254568156041
$ws.Range("F6").Value = 0.6960046614145909
$ws.Range("G6").Value = 0.5590941131201532
$ws.Range("H6").Value = 0.625176852670549
$ws.Range("I6").Value = 0.7080607263993954
$ws.Range("L6").Value = 0.1590095740485893
$ws.Range("O6").Value = 2.340951120064716
$ws.Range("C7").Value = 0.2218176138076871
$ws.Range("D7").Value = 0.006133347457781824
$ws.Range("E7").Value = 0.1345423215371859
$ws.Range("F7").Value = 0.7044598979951786
$ws.Range("G7").Value = 0.5681444924854873
$ws.Range("H7").Value = 0.6280763427505178
$ws.Range("I7").Value = 0.714764833313815
$ws.Range("L7").Value = 0.1603932375218164
$ws.Range("O7").Value = 2.366328777580151
$ws.Range("C8").Value = 0.2250378015571499
$ws.Range("D8").Value = 0.006263521966550556
$ws.Range("E8").Value = 0.1383194406000179
$ws.Range("F8").Value = 0.7427608516765076
$ws.Range("G8").Value = 0.6089053112631007
$ws.Range("H8").Value = 0.641623327560481
$ws.Range("I8").Value = 0.7453212505271409
$ws.Range("L8").Value = 0.1667203203117538
$ws.Range("O8").Value = 2.481602170323811
$ws.Range("C9").Value = 0.2322312903144024
$ws.Range("D9").Value = 0.006507425466285355
$ws.Range("E9").Value = 0.1463477039679475
$ws.Range("F9").Value = 0.821439386749816
$ws.Range("G9").Value = 0.6918155861135062
$ws.Range("H9").Value = 0.6708929033100901
$ws.Range("I9").Value = 0.8087480417135851
$ws.Range("L9").Value = 0.179926761494059
$ws.Range("O9").Value = 2.71950980883895
$ws.Range("C10").Value = 0.2380526478983285
$ws.Range("D10").Value = 0.006679668274886197
$ws.Range("E10").Value = 0.1526272775550055
$ws.Range("F10").Value = 0.8814422773570101
$ws.Range("G10").Value = 0.7545720582327533
$ws.Range("H10").Value = 0.6940505126202083
$ws.Range("I10").Value = 0.8575018524364708
$ws.Range("L10").Value = 0.1901217526640693
$ws.Range("O10").Value = 2.901595926179937
$ws.Range("C11").Value = 0.2408175807804298
$ws.Range("D11").Value = 0.006756512684301441
$ws.Range("E11").Value = 0.1555673673806339
$ws.Range("F11").Value = 0.9092213429262728
$ws.Range("G11").Value = 0.7835277731551002
$ws.Range("H11").Value = 0.7049457602524853
$ws.Range("I11").Value = 0.8801529968834814
$ws.Range("L11").Value = 0.1948677723273704
$ws.Range("O11").Value = 2.986031532157369
$ws.Range("C12").Value = 0.2418813775626916
$ws.Range("D12").Value = 0.006785393980866417
$ws.Range("E12").Value = 0.156692736421931
$ws.Range("F12").Value = 0.9198103212405471
$ws.Range("G12").Value = 0.7945514864978236
$ws.Range("H12").Value = 0.709123415987591
$ws.Range("I12").Value = 0.8887985810057728
$ws.Range("L12").Value = 0.196680603669094
$ws.Range("O12").Value = 3.018236523349401
$ws.Range("C13").Value = 0.2416515239705177
$ws.Range("D13").Value = 0.006779183587957149
$ws.Range("E13").Value = 0.1564498331702566
$ws.Range("F13").Value = 0.9175266955690518
$ws.Range("G13").Value = 0.792174713423293
$ws.Range("H13").Value = 0.7082213766728103
$ws.Range("I13").Value = 0.8869335699643983
$ws.Range("L13").Value = 0.1962894826811663
$ws.Range("O13").Value = 3.01129031868021
$ws.Range("C14").Value = 0.2409047638335409
$ws.Range("D14").Value = 0.006758893140791145
$ws.Range("E14").Value = 0.1556597111586768
$ws.Range("F14").Value = 0.9100911079304268
$ws.Range("G14").Value = 0.7844335207502411
$ws.Range("H14").Value = 0.7052884191755311
$ws.Range("I14").Value = 0.8808629091430902
$ws.Range("L14").Value = 0.195016601795146
$ws.Range("O14").Value = 2.988676422814763
$ws.Range("C15").Value = 0.2404495360860608
$ws.Range("D15").Value = 0.006746436237452968
$ws.Range("E15").Value = 0.1551773045707279
$ws.Range("F15").Value = 0.9055456700210698
$ws.Range("G15").Value = 0.7796994822955128
$ws.Range("H15").Value = 0.7034986516532626
$ws.Range("I15").Value = 0.8771533240768434
$ws.Range("L15").Value = 0.1942389606826538
$ws.Range("O15").Value = 2.974854862739676
$ws.Range("C16").Value = 0.2378743014690343
$ws.Range("D16").Value = 0.006674615853810906
$ws.Range("E16").Value = 0.1524368170123225
$ws.Range("F16").Value = 0.8796365894038018
$ws.Range("G16").Value = 0.7526879522886247
$ws.Range("H16").Value = 0.6933457404820444
$ws.Range("I16").Value = 0.8560310692201085
$ws.Range("L16").Value = 0.1898137718290513
$ws.Range("O16").Value = 2.896110169810186
$ws.Range("C17").Value = 0.2363243752458288
$ws.Range("D17").Value = 0.006630169038682965
$ws.Range("E17").Value = 0.1507770096439955
$ws.Range("F17").Value = 0.8638661227901139
$ws.Range("G17").Value = 0.7362217736789773
$ws.Range("H17").Value = 0.6872096428431007
$ws.Range("I17").Value = 0.8431944144987256
$ws.Range("L17").Value = 0.1871268257064287
$ws.Range("O17").Value = 2.848213919115551
$ws.Range("C18").Value = 0.2354438902956417
$ws.Range("D18").Value = 0.00660446235897183
$ws.Range("E18").Value = 0.1498301889176545
$ws.Range("F18").Value = 0.854840864300229
$ws.Range("G18").Value = 0.7267892042888775
$ws.Range("H18").Value = 0.6837142764616715
$ws.Range("I18").Value = 0.8358556055237187
$ws.Range("L18").Value = 0.1855915530932606
$ws.Range("O18").Value = 2.820816185126262
$ws.Range("C19").Value = 0.2351476614825145
$ws.Range("D19").Value = 0.006595734150334209
$ws.Range("E19").Value = 0.1495109605620044
$ws.Range("F19").Value = 0.8517928770304906
$ws.Range("G19").Value = 0.7236020768635001
$ws.Range("H19").Value = 0.6825366375807107
$ws.Range("I19").Value = 0.8333784499686203
$ws.Range("L19").Value = 0.1850734844566659
$ws.Range("O19").Value = 2.811565690694522
$ws.Range("C20").Value = 0.2364882300525721
$ws.Range("D20").Value = 0.006634915186035428
$ws.Range("E20").Value = 0.1509528857278468
$ws.Range("F20").Value = 0.8655402061964139
$ws.Range("G20").Value = 0.7379706567755591
$ws.Range("H20").Value = 0.6878593268105533
$ws.Range("I20").Value = 0.8445562924830909
$ws.Range("L20").Value = 0.1874118009699544
$ws.Range("O20").Value = 2.853296933605634
$ws.Range("C21").Value = 0.2411236501859406
$ws.Range("D21").Value = 0.006764858860632827
$ws.Range("E21").Value = 0.1558914627725514
$ws.Range("F21").Value = 0.9122732289703777
$ws.Range("G21").Value = 0.7867056987995511
$ws.Range("H21").Value = 0.7061484924561796
$ws.Range("I21").Value = 0.882644160706235
$ws.Range("L21").Value = 0.195390053550085
$ws.Range("O21").Value = 2.995312397728412
$ws.Range("C22").Value = 0.2442509461550486
$ws.Range("D22").Value = 0.006848513406062295
$ws.Range("E22").Value = 0.1591891827106195
$ws.Range("F22").Value = 0.9432219894342353
$ws.Range("G22").Value = 0.8188997715168966
$ws.Range("H22").Value = 0.718403805864483
$ws.Range("I22").Value = 0.9079336965696569
$ws.Range("L22").Value = 0.200695351301718
$ws.Range("O22").Value = 3.089474715421375
$ws.Range("C23").Value = 0.2425729071831881
$ws.Range("D23").Value = 0.006803982008655396
$ws.Range("E23").Value = 0.1574227102323107
$ws.Range("F23").Value = 0.9266668589498721
$ws.Range("G23").Value = 0.8016857412672209
$ws.Range("H23").Value = 0.7118352587157517
$ws.Range("I23").Value = 0.8943998444289889
$ws.Range("L23").Value = 0.1978554671858177
$ws.Range("O23").Value = 3.039095115121427
$ws.Range("C24").Value = 0.2364141183130499
$ws.Range("D24").Value = 0.006632769931339766
$ws.Range("E24").Value = 0.1508733490257086
$ws.Range("F24").Value = 0.8647832241590692
$ws.Range("G24").Value = 0.737179880661273
$ws.Range("H24").Value = 0.6875655038883508
$ws.Range("I24").Value = 0.843940459256288
$ws.Range("L24").Value = 0.187282934103294
$ws.Range("O24").Value = 2.850998471484786
$ws.Range("C25").Value = 0.2301911261090197
$ws.Range("D25").Value = 0.006442662609300598
$ws.Range("E25").Value = 0.1441090771027689
$ws.Range("F25").Value = 0.7997703333636537
$ws.Range("G25").Value = 0.6690645539529498
$ws.Range("H25").Value = 0.6626847321016385
$ws.Range("I25").Value = 0.7912123505931419
$ws.Range("L25").Value = 0.1762679991570337
$ws.Range("O25").Value = 2.653873203053593
